$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,7).Value = 0.009318
$ws.Cells.Item(2,9).Value = 0.1755772177976535
$ws.Cells.Item(2,10).Value = 0.1755772177976535
$ws.Cells.Item(2,13).Value = 33.041958
$ws.Cells.Item(2,14).Value = 99.125874
$ws.Cells.Item(2,15).Value = 0.8362572961276681
$ws.Cells.Item(2,16).Value = 0.8362572961276681
$ws.Cells.Item(2,17).Value = 0.307884964644
$ws.Cells.Item(2,18).Value = 2.770964681796
$ws.Cells.Item(2,19).Value = 0.1468277294170844
$ws.Cells.Item(2,20).Value = 0.1468277294170844

# Row 3
$ws.Cells.Item(3,7).Value = 0.009318
$ws.Cells.Item(3,9).Value = 0.1755772177976535
$ws.Cells.Item(3,10).Value = 0.1755772177976535
$ws.Cells.Item(3,13).Value = 2.582190666666666
$ws.Cells.Item(3,14).Value = 7.746571999999999
$ws.Cells.Item(3,15).Value = 0.06535253706795362
$ws.Cells.Item(3,16).Value = 0.06535253706795363
$ws.Cells.Item(3,17).Value = 0.02406085263199999
$ws.Cells.Item(3,18).Value = 0.216547673688
$ws.Cells.Item(3,19).Value = 0.01147441663440931
$ws.Cells.Item(3,20).Value = 0.01147441663440931

# Row 4
$ws.Cells.Item(4,7).Value = 0.009318
$ws.Cells.Item(4,9).Value = 0.1755772177976535
$ws.Cells.Item(4,10).Value = 0.1755772177976535
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.2536366666666667
$ws.Cells.Item(4,14).Value = 0.76091
$ws.Cells.Item(4,15).Value = 0.006419277969710552
$ws.Cells.Item(4,16).Value = 0.006419277969710552
$ws.Cells.Item(4,17).Value = 0.00236338646
$ws.Cells.Item(4,18).Value = 0.02127047814
$ws.Cells.Item(4,19).Value = 0.001127078966191548
$ws.Cells.Item(4,20).Value = 0.001127078966191548

# Row 5
$ws.Cells.Item(5,7).Value = 0.009318
$ws.Cells.Item(5,9).Value = 0.1755772177976535
$ws.Cells.Item(5,10).Value = 0.1755772177976535
$ws.Cells.Item(5,13).Value = 1.932675333333333
$ws.Cells.Item(5,14).Value = 5.798026
$ws.Cells.Item(5,15).Value = 0.04891398531969482
$ws.Cells.Item(5,16).Value = 0.04891398531969483
$ws.Cells.Item(5,17).Value = 0.018008668756
$ws.Cells.Item(5,18).Value = 0.162078018804
$ws.Cells.Item(5,19).Value = 0.008588181453827282
$ws.Cells.Item(5,20).Value = 0.008588181453827282

# Row 6
$ws.Cells.Item(6,7).Value = 0.009318
$ws.Cells.Item(6,9).Value = 0.1755772177976535
$ws.Cells.Item(6,10).Value = 0.1755772177976535
$ws.Cells.Item(6,13).Value = 1.701252
$ws.Cells.Item(6,14).Value = 5.103756
$ws.Cells.Item(6,15).Value = 0.04305690351497292
$ws.Cells.Item(6,16).Value = 0.04305690351497292
$ws.Cells.Item(6,17).Value = 0.015852266136
$ws.Cells.Item(6,18).Value = 0.142670395224
$ws.Cells.Item(6,19).Value = 0.007559811326140951
$ws.Cells.Item(6,20).Value = 0.007559811326140951

# Row 7
$ws.Cells.Item(7,7).Value = 0.03948166666666666
$ws.Cells.Item(7,8).Value = 0.118445
$ws.Cells.Item(7,9).Value = 0.743945180011557
$ws.Cells.Item(7,10).Value = 0.743945180011557
$ws.Cells.Item(7,13).Value = 33.041958
$ws.Cells.Item(7,14).Value = 99.125874
$ws.Cells.Item(7,15).Value = 0.8362572961276681
$ws.Cells.Item(7,16).Value = 0.8362572961276681
$ws.Cells.Item(7,17).Value = 1.30455157177
$ws.Cells.Item(7,18).Value = 11.74096414593
$ws.Cells.Item(7,19).Value = 0.622129584703676
$ws.Cells.Item(7,20).Value = 0.622129584703676

# Row 8
$ws.Cells.Item(8,7).Value = 0.03948166666666666
$ws.Cells.Item(8,8).Value = 0.118445
$ws.Cells.Item(8,9).Value = 0.743945180011557
$ws.Cells.Item(8,10).Value = 0.743945180011557
$ws.Cells.Item(8,13).Value = 2.582190666666666
$ws.Cells.Item(8,14).Value = 7.746571999999999
$ws.Cells.Item(8,15).Value = 0.06535253706795362
$ws.Cells.Item(8,16).Value = 0.06535253706795363
$ws.Cells.Item(8,17).Value = 0.1019491911711111
$ws.Cells.Item(8,18).Value = 0.9175427205399999
$ws.Cells.Item(8,19).Value = 0.0486187049532307
$ws.Cells.Item(8,20).Value = 0.04861870495323072

# Row 9
$ws.Cells.Item(9,7).Value = 0.03948166666666666
$ws.Cells.Item(9,8).Value = 0.118445
$ws.Cells.Item(9,9).Value = 0.743945180011557
$ws.Cells.Item(9,10).Value = 0.743945180011557
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.2536366666666667
$ws.Cells.Item(9,14).Value = 0.76091
$ws.Cells.Item(9,15).Value = 0.006419277969710552
$ws.Cells.Item(9,16).Value = 0.006419277969710552
$ws.Cells.Item(9,17).Value = 0.01001399832777778
$ws.Cells.Item(9,18).Value = 0.09012598494999999
$ws.Cells.Item(9,19).Value = 0.004775590904720538
$ws.Cells.Item(9,20).Value = 0.004775590904720538

# Row 10
$ws.Cells.Item(10,7).Value = 0.03948166666666666
$ws.Cells.Item(10,8).Value = 0.118445
$ws.Cells.Item(10,9).Value = 0.743945180011557
$ws.Cells.Item(10,10).Value = 0.743945180011557
$ws.Cells.Item(10,13).Value = 1.932675333333333
$ws.Cells.Item(10,14).Value = 5.798026
$ws.Cells.Item(10,15).Value = 0.04891398531969482
$ws.Cells.Item(10,16).Value = 0.04891398531969483
$ws.Cells.Item(10,17).Value = 0.07630524328555555
$ws.Cells.Item(10,18).Value = 0.68674718957
$ws.Cells.Item(10,19).Value = 0.03638932361374302
$ws.Cells.Item(10,20).Value = 0.03638932361374302

# Row 11
$ws.Cells.Item(11,7).Value = 0.03948166666666666
$ws.Cells.Item(11,8).Value = 0.118445
$ws.Cells.Item(11,9).Value = 0.743945180011557
$ws.Cells.Item(11,10).Value = 0.743945180011557
$ws.Cells.Item(11,13).Value = 1.701252
$ws.Cells.Item(11,14).Value = 5.103756
$ws.Cells.Item(11,15).Value = 0.04305690351497292
$ws.Cells.Item(11,16).Value = 0.04305690351497292
$ws.Cells.Item(11,17).Value = 0.06716826438
$ws.Cells.Item(11,18).Value = 0.60451437942
$ws.Cells.Item(11,19).Value = 0.03203197583618677
$ws.Cells.Item(11,20).Value = 0.03203197583618677

# Row 12
$ws.Cells.Item(12,7).Value = 0.004271
$ws.Cells.Item(12,8).Value = 0.012813
$ws.Cells.Item(12,9).Value = 0.08047760219078964
$ws.Cells.Item(12,10).Value = 0.08047760219078964
$ws.Cells.Item(12,13).Value = 33.041958
$ws.Cells.Item(12,14).Value = 99.125874
$ws.Cells.Item(12,15).Value = 0.8362572961276681
$ws.Cells.Item(12,16).Value = 0.8362572961276681
$ws.Cells.Item(12,17).Value = 0.141122202618
$ws.Cells.Item(12,18).Value = 1.270099823562
$ws.Cells.Item(12,19).Value = 0.06729998200690784
$ws.Cells.Item(12,20).Value = 0.06729998200690784

# Row 13
$ws.Cells.Item(13,7).Value = 0.004271
$ws.Cells.Item(13,8).Value = 0.012813
$ws.Cells.Item(13,9).Value = 0.08047760219078964
$ws.Cells.Item(13,10).Value = 0.08047760219078964
$ws.Cells.Item(13,13).Value = 2.582190666666666
$ws.Cells.Item(13,14).Value = 7.746571999999999
$ws.Cells.Item(13,15).Value = 0.06535253706795362
$ws.Cells.Item(13,16).Value = 0.06535253706795363
$ws.Cells.Item(13,17).Value = 0.01102853633733333
$ws.Cells.Item(13,18).Value = 0.09925682703599999
$ws.Cells.Item(13,19).Value = 0.005259415480313606
$ws.Cells.Item(13,20).Value = 0.005259415480313607

# Row 14
$ws.Cells.Item(14,7).Value = 0.004271
$ws.Cells.Item(14,8).Value = 0.012813
$ws.Cells.Item(14,9).Value = 0.08047760219078964
$ws.Cells.Item(14,10).Value = 0.08047760219078964
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 0.2536366666666667
$ws.Cells.Item(14,14).Value = 0.76091
$ws.Cells.Item(14,15).Value = 0.006419277969710552
$ws.Cells.Item(14,16).Value = 0.006419277969710552
$ws.Cells.Item(14,17).Value = 0.001083282203333333
$ws.Cells.Item(14,18).Value = 0.009749539829999999
$ws.Cells.Item(14,19).Value = 0.0005166080987984656
$ws.Cells.Item(14,20).Value = 0.0005166080987984656

# Row 15
$ws.Cells.Item(15,7).Value = 0.004271
$ws.Cells.Item(15,8).Value = 0.012813
$ws.Cells.Item(15,9).Value = 0.08047760219078964
$ws.Cells.Item(15,10).Value = 0.08047760219078964
$ws.Cells.Item(15,13).Value = 1.932675333333333
$ws.Cells.Item(15,14).Value = 5.798026
$ws.Cells.Item(15,15).Value = 0.04891398531969482
$ws.Cells.Item(15,16).Value = 0.04891398531969483
$ws.Cells.Item(15,17).Value = 0.008254456348666665
$ws.Cells.Item(15,18).Value = 0.074290107138
$ws.Cells.Item(15,19).Value = 0.003936480252124524
$ws.Cells.Item(15,20).Value = 0.003936480252124524

# Row 16
$ws.Cells.Item(16,7).Value = 0.004271
$ws.Cells.Item(16,8).Value = 0.012813
$ws.Cells.Item(16,9).Value = 0.08047760219078964
$ws.Cells.Item(16,10).Value = 0.08047760219078964
$ws.Cells.Item(16,13).Value = 1.701252
$ws.Cells.Item(16,14).Value = 5.103756
$ws.Cells.Item(16,15).Value = 0.04305690351497292
$ws.Cells.Item(16,16).Value = 0.04305690351497292
$ws.Cells.Item(16,17).Value = 0.007266047291999999
$ws.Cells.Item(16,18).Value = 0.06539442562799999
$ws.Cells.Item(16,19).Value = 0.003465116352645203
$ws.Cells.Item(16,20).Value = 0.003465116352645203

